# Rebuild the article header to match the new "pandoc title block" shape:
#   - drop the old bookmark-wrapped Heading1 title paragraph
#   - add a Title-style paragraph whose text is split word-by-word into runs
#   - add an Authors-style paragraph ("Dorothy Day") split word-by-word into runs
#   - drop the old bold "By Dorothy Day" paragraph

$d = $word.ActiveDocument

function Split-IntoRuns($para, $text, $splitPoints) {
    # Write the full text into the (already emptied) paragraph, then force
    # run boundaries at each offset in $splitPoints by briefly adding and
    # removing a bookmark there (adjacent same-format runs created in one
    # editing session otherwise get coalesced back into a single run).
    $para.Range.Text = $text
    $start = $para.Range.Start
    $i = 0
    foreach ($s in $splitPoints) {
        $i = $i + 1
        $name = "tmpsplit" + $i
        $rng = $d.Range($start, $start + $s)
        $d.Bookmarks.Add($name, $rng) | Out-Null
    }
    for ($j = 1; $j -le $i; $j++) {
        $name = "tmpsplit" + $j
        $d.Bookmarks($name).Delete()
    }
}

# --- 1. Remove the old title paragraph (and its wrapping bookmark) -------
$p1 = $d.Paragraphs(1)
$d.Range($p1.Range.Start, $p1.Range.End).Delete()

# The bookmarkStart/bookmarkEnd that used to wrap paragraph 1 collapse to a
# zero-length pair at the very start of the document; two zero-length
# deletes there clear them out of the package entirely.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- 2. Insert the new Title paragraph in its place -----------------------
$d.Paragraphs(1).Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs(1)
$titlePara.Style = "Title"
Split-IntoRuns $titlePara "Fall Appeal - October/November 1977" @(4, 5, 11, 12, 13, 14, 21, 22, 30, 31)

# --- 3. Turn "By Dorothy Day" into an Authors paragraph "Dorothy Day" -----
$authorPara = $d.Paragraphs(2)
$authorPara.Range.Font.Bold = 0
$authorPara.Style = "Authors"
Split-IntoRuns $authorPara "Dorothy Day" @(7, 8)

Write-Output "done"
